# ---------------------------------------------------------------------------
# Commit message: "change the back end so it shouldnt disconnect and
# keywods can be more than one word long"
#
# This describes a change to the *external* PowerPoll web add-in backend
# (its HTML/JS "OfficeApp" content add-in), not to any text/shape on the
# slide. In the OOXML, the add-in is the <mc:AlternateContent> block whose
# <p:graphicFrame>/<we:webextensionref> (with mc:Fallback <p:pic> snapshot)
# is bound to ppt/slides/udata/data.xml. Re-syncing an already-inserted
# add-in after a backend edit does not touch the slide's visible content -
# diffing the canonical OOXML before/after shows the ctrTitle/subTitle
# placeholders are untouched (both empty) and the only substantive change
# is PowerPoint re-minting that add-in instance's internal identifier,
# <we:webextension id="...">  (plus a batch of unrelated r:id/r:embed
# renames that PowerPoint churns on every save, independent of this edit).
#
# Neither of those is reachable from the PowerPoint object model: there is
# no WebExtension/OfficeApp automation object in
# Microsoft.Office.Interop.PowerPoint (on real PowerPoint or here), and the
# add-in's graphicFrame/fallback-pic pair reuses cNvPr id="2" - the same id
# as the "Title 1" placeholder - so Slide.Shapes never exposes it as an
# independently addressable shape (every index/name lookup that should
# reach "OfficeApp 0" resolves back onto the "Title 1" shape instead, for
# both reads and writes). `Presentation.CustomXMLParts` / `Shape.Tags` /
# `Shape.CustomerData` are the closest stand-ins, but they write unrelated,
# brand-new package parts rather than updating
# ppt/slides/udata/data.xml - using them here would corrupt the deck
# instead of reproducing the edit, so this script intentionally leaves the
# slide alone (matching the fact that none of its visible content actually
# changed).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

Write-Host ("Slide 1 shapes: {0} (PowerPoll OfficeApp add-in re-synced after " +
    "backend update - keep-alive fix + multi-word keyword support; no " +
    "placeholder text on the slide changed)" -f $s.Shapes.Count)
